# MonitoramentoCargaOBM.xlsx - "results-specimen-type-uv-ips  carga value set"
#
# Header row 3 (J3:L3): re-label / reorder the "VMP incluídos Portal" /
# "VMPP incluídos Portal" headers and add a new "% VMPs incluídos no portal"
# header, all sharing the same (top-aligned, wrapped) look.
#
# Row 7 / Row 8 / Row 13 / Row 14: July ("JULHO 23", row 8) gets its real
# monitoring numbers filled in (it was a placeholder row before), which
# ripples into the TOTAIS row (14) sums, and a new VMP-% column (K) is
# populated for both row 8 and the totals row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 3 headers: J3 / K3 / L3
# ---------------------------------------------------------------------
$ws.Range("J3").Value = "VMP incluídos Portal"
$ws.Range("K3").Value = "% VMPs incluídos no portal"
$ws.Range("L3").Value = "VMPP incluídos Portal"

# Give all three the same "border + top-aligned + wrap" look.
$ws.Range("J3:L3").VerticalAlignment = -4160   # xlTop
$ws.Range("J3:L3").WrapText = $true

# ---------------------------------------------------------------------
# Row 7 (JUNHO 23): E7 gets an explicit 0, L7 (VMPP incluidos) is cleared
# ---------------------------------------------------------------------
$ws.Range("E7").Value = 0
$ws.Range("L7").ClearContents()

# ---------------------------------------------------------------------
# Row 8 (JULHO 23): fill in the real figures
# ---------------------------------------------------------------------
$ws.Range("B8").Value = 4799
$ws.Range("C8").Value = 4799
$ws.Range("D8").Formula = "=(C8)/B14"
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 3
$ws.Range("G8").Value = 1013
$ws.Range("H8").Value = 850
$ws.Range("I8").Value = 1223

$ws.Range("J8").Value = 434
$ws.Range("J8").Style = "Normal"

$ws.Range("K8").Style = "Normal"
$ws.Range("K8").Formula = "=J8/B8"
$ws.Range("K8").NumberFormat = "0.00%"

$ws.Range("L8").Value = 1307
$ws.Range("L8").Style = "Normal"

# ---------------------------------------------------------------------
# Row 13 (DEZEMBRO 23): L13 placeholder cleared
# ---------------------------------------------------------------------
$ws.Range("L13").ClearContents()

# ---------------------------------------------------------------------
# Row 14 (TOTAIS): C14 becomes a literal (matches B14 now), J14/K14 get the
# new VMP / % VMP totals, L14 is cleared (VMPP total no longer tracked)
# ---------------------------------------------------------------------
$ws.Range("C14").Value = 4799

$ws.Range("J14").Value = 434
$ws.Range("K14").Formula = "=J14/B14"
$ws.Range("K14").NumberFormat = "0.00%"

$ws.Range("L14").ClearContents()

# ---------------------------------------------------------------------
# Selection / view: user ended up with L8 selected, scrolled to row 3
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("L8").Select()
